# repull data, push all data, mean calculation
# Update dSF column (F) values for several rows based on repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -2
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = 1
